$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.959.02'
$ws.Range('E2').Value = '  +2.22%  '
$ws.Range('D3').Value = '3.185.07'
$ws.Range('E3').Value = '  +0.97%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '535.60'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.22%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.95'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.03%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.23%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.530'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.13%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.34'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.39%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.112'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.00%  '
$ws.Range('E11').Value = '  -1.52%  '
$ws.Range('D12').Value = '3.734.70'
$ws.Range('E12').Value = '  +0.90%  '
$ws.Range('E13').Value = '  -2.49%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.76'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.14%  '
$ws.Range('E15').Value = '  +0.11%  '
$ws.Range('D16').Value = '59.983.57'
$ws.Range('E16').Value = '  +2.15%  '
$ws.Range('D17').Value = '3.196.44'
$ws.Range('E17').Value = '  +1.76%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.24'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.00%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.23'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '8.17'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.56%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '368.84'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.75%  '
$ws.Range('E22').Value = '  +0.02%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.522'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.46%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '69.35'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.31%  '
$ws.Range('E25').Value = '  +1.03%  '
$ws.Range('B26').Value = 'Binance-PegBSC-USD'
$ws.Range('C26').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.15%  '
$ws.Range('B27').Value = 'InternetComputer(DFINITY)'
$ws.Range('C27').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.53'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.95%  '
$ws.Range('D28').Value = '0.0₃0871'
$ws.Range('E28').Value = '  +1.00%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '22.45'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.85%  '
$ws.Range('E30').Value = '  +0.51%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.09'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.28%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.26'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.62%  '
$ws.Range('B33').Value = 'Aptos'
$ws.Range('C33').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.58'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.51%  '
$ws.Range('B34').Value = 'Fetch.AI'
$ws.Range('C34').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.19'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.55%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '157.56'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.44%  '
$ws.Range('E36').Value = '  +1.66%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '26.42'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +5.39%  '
$ws.Range('D38').Value = '2.783.56'
$ws.Range('E38').Value = '  +5.99%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0705'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.91%  '
$ws.Range('E40').Value = '  +7.81%  '
$ws.Range('E41').Value = '  +0.53%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.23'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.29%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '39.72'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.71%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.718'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.21%  '
$ws.Range('E45').Value = '  +0.75%  '
$ws.Range('D46').Value = '3.227.51'
$ws.Range('E46').Value = '  +0.94%  '
$ws.Range('E47').Value = '  +0.39%  '
$ws.Range('E48').Value = '  -0.81%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '20.59'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.18%  '
$ws.Range('E50').Value = '  +5.42%  '
